# Update the user-entered value for tau (B5) on Sheet1 from 1 to 0.8.
# Dependent formulas (A2, B2, B6, B10) recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = 0.8

$excel.CalculateFullRebuild()
